$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Schedule")
$ws2 = $wb.Worksheets.Item("Detailed")

$ws1.Range("E3").Value = 819.03405675
$ws1.Range("F3").Value = 48.15015030864198
$ws1.Range("E4").Value = 1294.16421225
$ws1.Range("F4").Value = 45.64953129629629

$ws2.Range("B37").Value = 87.84635
$ws2.Range("B38").Value = 85.01606
$ws2.Range("B39").Value = 100.01
$ws2.Range("B40").Value = 108.89
$ws2.Range("C40").Value = "historical"
$ws2.Range("B41").Value = 119.51613
$ws2.Range("C41").Value = "historical"
$ws2.Range("B42").Value = 108.89
$ws2.Range("C42").Value = "historical"
$ws2.Range("B43").Value = 120.10153
$ws2.Range("C43").Value = "historical"
$ws2.Range("B44").Value = 108.89
$ws2.Range("C44").Value = "historical"
$ws2.Range("B45").Value = 105.79
$ws2.Range("C45").Value = "historical"
$ws2.Range("B46").Value = 86.07659
$ws2.Range("C46").Value = "historical"
$ws2.Range("C47").Value = "historical"
$ws2.Range("B48").Value = 108.01
$ws2.Range("C48").Value = "historical"
$ws2.Range("B49").Value = 104.83796
$ws2.Range("B50").Value = 85.65000000000001
$ws2.Range("B51").Value = 105.79
$ws2.Range("B52").Value = 105.79
$ws2.Range("B53").Value = 95.00112
$ws2.Range("B54").Value = 85.65000000000001
$ws2.Range("B55").Value = 85.65000000000001
$ws2.Range("B56").Value = 87.01595
$ws2.Range("B57").Value = 85.65000000000001
$ws2.Range("B58").Value = 103.83786
$ws2.Range("B59").Value = 108.01
$ws2.Range("B61").Value = 131.52376
$ws2.Range("B62").Value = 130.19921
$ws2.Range("B63").Value = 149.71495
$ws2.Range("B64").Value = 147.51
$ws2.Range("B65").Value = 147.51
$ws2.Range("B66").Value = 138
$ws2.Range("B67").Value = 107.88115
$ws2.Range("B68").Value = 100.01
$ws2.Range("B69").Value = 108.89
$ws2.Range("B70").Value = 105.79
$ws2.Range("B71").Value = 85.65000000000001
$ws2.Range("B73").Value = 78.0001
$ws2.Range("B74").Value = 79.44887
$ws2.Range("B78").Value = 84.33557999999999
$ws2.Range("B80").Value = 82.08329000000001
$ws2.Range("B81").Value = 98.46892
$ws2.Range("B82").Value = 100.25304
$ws2.Range("B83").Value = 93.76103999999999
$ws2.Range("B85").Value = 107.62964
$ws2.Range("B86").Value = 114.75372
$ws2.Range("B87").Value = 130.86899
$ws2.Range("B88").Value = 166.99
$ws2.Range("B90").Value = 173.3557
$ws2.Range("B91").Value = 149.06831
$ws2.Range("B93").Value = 147.51
$ws2.Range("B94").Value = 111.07152
$ws2.Range("B95").Value = 136.74568
$ws2.Range("B96").Value = 108.89
$ws2.Range("B97").Value = 105.79

Write-Host "Applied run 164 updates"
